$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 16671146
$ws.Range("I62").Value = 25004120
$ws.Range("K62").Value = 25004120
$ws.Range("M62").Value = -25003496
$ws.Range("H65").Value = 16671146
$ws.Range("I65").Value = 25004120
$ws.Range("K65").Value = 125020600
$ws.Range("M65").Value = -125017480
$ws.Range("H74").Value = 5589.9575
$ws.Range("I74").Value = 5022.4346
$ws.Range("K74").Value = 5022.4346
$ws.Range("M74").Value = -4086.4346
$ws.Range("H77").Value = 5589.9575
$ws.Range("I77").Value = 5022.4346
$ws.Range("K77").Value = 25112.173
$ws.Range("M77").Value = -20432.173
$ws.Range("H82").Value = 8249.5
$ws.Range("J82").Value = 20998
$ws.Range("L82").Value = 62994
$ws.Range("N82").Value = -63806
$ws.Range("H85").Value = 8249.5
$ws.Range("J85").Value = 20998
$ws.Range("L85").Value = 62994
$ws.Range("N85").Value = -65802
$ws.Range("H86").Value = 4474.6
$ws.Range("I86").Value = 5332.6665
$ws.Range("J86").Value = 4106.857
$ws.Range("K86").Value = 5332.6665
$ws.Range("L86").Value = 4106.857
$ws.Range("M86").Value = -4209.6665
$ws.Range("N86").Value = -6352.857
$ws.Range("H89").Value = 4474.6
$ws.Range("I89").Value = 5332.6665
$ws.Range("J89").Value = 4106.857
$ws.Range("K89").Value = 26663.3325
$ws.Range("L89").Value = 20534.285
$ws.Range("M89").Value = -21047.3325
$ws.Range("N89").Value = -31766.285
$ws.Range("H106").Value = 10110554
$ws.Range("I106").Value = 18195012
$ws.Range("K106").Value = 18195012
$ws.Range("M106").Value = -18194381
$ws.Range("H113").Value = 5987.355
$ws.Range("I113").Value = 5517.381
$ws.Range("J113").Value = 6974.3
$ws.Range("K113").Value = 5517.381
$ws.Range("L113").Value = 6974.3
$ws.Range("M113").Value = -2263.381
$ws.Range("N113").Value = -13482.3
$ws.Range("H116").Value = 7404.4287
$ws.Range("I116").Value = 8407.111000000001
$ws.Range("K116").Value = 8407.111000000001
$ws.Range("M116").Value = -4965.111000000001
$ws.Range("H141").Value = 5476.5864
$ws.Range("I141").Value = 4232.96
$ws.Range("K141").Value = 12698.88
$ws.Range("M141").Value = -7518.880000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 881.6
$ws.Range("J20").Value = 822.8333
$ws.Range("L20").Value = 822.8333
$ws.Range("N20").Value = -1316.8333
$ws.Range("H107").Value = 1767.2
$ws.Range("I107").Value = 1511.2858
$ws.Range("J107").Value = 2364.3333
$ws.Range("K107").Value = 1511.2858
$ws.Range("L107").Value = 2364.3333
$ws.Range("M107").Value = 408.7141999999999
$ws.Range("N107").Value = -6204.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I94").Value = 2128.4443
$ws.Range("J94").Value = 2159.4614
$ws.Range("K94").Value = 2128.4443
$ws.Range("L94").Value = 2159.4614
$ws.Range("M94").Value = -1677.4443
$ws.Range("N94").Value = -3061.4614
$ws.Range("H105").Value = 1537.4
$ws.Range("I105").Value = 899
$ws.Range("K105").Value = 899
$ws.Range("M105").Value = 848
$ws.Range("H107").Value = 393.8421
$ws.Range("I107").Value = 257.05884
$ws.Range("J107").Value = 1556.5
$ws.Range("K107").Value = 257.05884
$ws.Range("L107").Value = 1556.5
$ws.Range("M107").Value = 1662.94116
$ws.Range("N107").Value = -5396.5
$ws.Range("H132").Value = 3867
$ws.Range("I132").Value = 3687.3044
$ws.Range("K132").Value = 11061.9132
$ws.Range("M132").Value = -8531.913199999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 171.72
$ws.Range("I2").Value = 35.411766
$ws.Range("J2").Value = 461.375
$ws.Range("K2").Value = 212.470596
$ws.Range("L2").Value = 2768.25
$ws.Range("M2").Value = -99.470596
$ws.Range("N2").Value = -2994.25
$ws.Range("H23").Value = 915.4
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 915.4
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 2746.2
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -3216.2
$ws.Range("H38").Value = 317.05264
$ws.Range("I38").Value = 258.18182
$ws.Range("J38").Value = 398
$ws.Range("K38").Value = 774.54546
$ws.Range("L38").Value = 1194
$ws.Range("M38").Value = -427.54546
$ws.Range("N38").Value = -1888
$ws.Range("H80").Value = 2672.125
$ws.Range("I80").Value = 1999.5
$ws.Range("J80").Value = 2896.3333
$ws.Range("K80").Value = 5998.5
$ws.Range("L80").Value = 8688.999899999999
$ws.Range("M80").Value = -5062.5
$ws.Range("N80").Value = -10560.9999
$ws.Range("H83").Value = 2672.125
$ws.Range("I83").Value = 1999.5
$ws.Range("J83").Value = 2896.3333
$ws.Range("K83").Value = 17995.5
$ws.Range("L83").Value = 26066.9997
$ws.Range("M83").Value = -13315.5
$ws.Range("N83").Value = -35426.9997
$ws.Range("H86").Value = 7226.6
$ws.Range("I86").Value = 600
$ws.Range("K86").Value = 1800
$ws.Range("M86").Value = -614
$ws.Range("H89").Value = 7226.6
$ws.Range("I89").Value = 600
$ws.Range("K89").Value = 5400
$ws.Range("M89").Value = 528
$ws.Range("H107").Value = 630.8214
$ws.Range("J107").Value = 640.9231
$ws.Range("L107").Value = 1922.7693
$ws.Range("N107").Value = -5762.7693
$ws.Range("H132").Value = 1262.303
$ws.Range("J132").Value = 1665.3125
$ws.Range("L132").Value = 14987.8125
$ws.Range("N132").Value = -20047.8125
$ws.Range("H133").Value = 21620.318
$ws.Range("I133").Value = 3663.8572
$ws.Range("K133").Value = 10991.5716
$ws.Range("M133").Value = -5931.571599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10820.9
$ws.Range("J70").Value = 5813.5713
$ws.Range("L70").Value = 5813.5713
$ws.Range("N70").Value = -6353.5713
$ws.Range("H73").Value = 10820.9
$ws.Range("J73").Value = 5813.5713
$ws.Range("L73").Value = 5813.5713
$ws.Range("N73").Value = -7685.5713
$ws.Range("H92").Value = 8705.833000000001
$ws.Range("J92").Value = 8705.833000000001
$ws.Range("L92").Value = 8705.833000000001
$ws.Range("N92").Value = -12449.833
$ws.Range("H107").Value = 499.5
$ws.Range("I107").Value = 499.5
$ws.Range("K107").Value = 499.5
$ws.Range("M107").Value = 1420.5
$ws.Range("H113").Value = 1644.3158
$ws.Range("I113").Value = 1030.7693
$ws.Range("K113").Value = 1030.7693
$ws.Range("M113").Value = 1139.2307
$ws.Range("H126").Value = 5677.3887
$ws.Range("J126").Value = 5604.9
$ws.Range("L126").Value = 16814.7
$ws.Range("N126").Value = -21754.7
$ws.Range("H132").Value = 2305.7354
$ws.Range("I132").Value = 1763.9
$ws.Range("K132").Value = 5291.700000000001
$ws.Range("M132").Value = -2761.700000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1112.375
$ws.Range("J22").Value = 1436.3334
$ws.Range("L22").Value = 1436.3334
$ws.Range("N22").Value = -2026.3334
$ws.Range("H27").Value = 1112.375
$ws.Range("J27").Value = 1436.3334
$ws.Range("L27").Value = 1436.3334
$ws.Range("N27").Value = -1650.3334
$ws.Range("H46").Value = 1946.0667
$ws.Range("I46").Value = 1211.1111
$ws.Range("K46").Value = 1211.1111
$ws.Range("M46").Value = -1023.1111
$ws.Range("H55").Value = 789.03705
$ws.Range("I55").Value = 141.76923
$ws.Range("K55").Value = 141.76923
$ws.Range("M55").Value = 31.23077000000001
$ws.Range("H132").Value = 3225.4614
$ws.Range("I132").Value = 2466.6924
$ws.Range("J132").Value = 3984.2307
$ws.Range("K132").Value = 7400.0772
$ws.Range("L132").Value = 11952.6921
$ws.Range("M132").Value = -4870.0772
$ws.Range("N132").Value = -17012.6921
$ws.Range("H136").Value = 2236.8438
$ws.Range("I136").Value = 2063.4614
$ws.Range("K136").Value = 6190.3842
$ws.Range("M136").Value = -3640.3842

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 32250
$ws.Range("J103").Value = 32250
$ws.Range("L103").Value = 32250
$ws.Range("N103").Value = -34594
$ws.Range("H122").Value = 8219.474
$ws.Range("I122").Value = 5881.853
$ws.Range("K122").Value = 17645.559
$ws.Range("M122").Value = -15195.559
$ws.Range("H132").Value = 2712.5
$ws.Range("I132").Value = 2712.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8137.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5607.5
$ws.Range("N132").ClearContents()
